# Auto-generated edit script applying the Hyperion_Profits market-data refresh.
# For each changed cell: set new value; two rows (ARM!53, ARM!55) also
# gain/lose a profit cell entirely (M53 added, N55 removed), matching the
# source diff's added/removed <c> elements.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 863.8409
$ws.Range("I15").Value = 863.8409
$ws.Range("K15").Value = 2591.5227
$ws.Range("M15").Value = -2422.5227
$ws.Range("H19").Value = 1575
$ws.Range("I19").Value = 250
$ws.Range("J19").Value = 1953.5714
$ws.Range("K19").Value = 250
$ws.Range("L19").Value = 1953.5714
$ws.Range("M19").Value = -75
$ws.Range("N19").Value = -2303.5714
$ws.Range("H62").Value = 6520.8667
$ws.Range("J62").Value = 7193.3076
$ws.Range("L62").Value = 7193.3076
$ws.Range("N62").Value = -8441.3076
$ws.Range("H65").Value = 6520.8667
$ws.Range("J65").Value = 7193.3076
$ws.Range("L65").Value = 35966.538
$ws.Range("N65").Value = -42206.538
$ws.Range("H112").Value = 5287.7144
$ws.Range("J112").Value = 6232.913
$ws.Range("L112").Value = 18698.739
$ws.Range("N112").Value = -20914.739
$ws.Range("H113").Value = 5299.3667
$ws.Range("I113").Value = 3595.5557
$ws.Range("K113").Value = 3595.5557
$ws.Range("M113").Value = -341.5556999999999
$ws.Range("H115").Value = 1509.3334
$ws.Range("I115").Value = 1004.5455
$ws.Range("K115").Value = 3013.6365
$ws.Range("M115").Value = -1446.6365
$ws.Range("H130").Value = 76520
$ws.Range("J130").Value = 76520
$ws.Range("L130").Value = 76520
$ws.Range("N130").Value = -86560
$ws.Range("H135").Value = 814.90625
$ws.Range("I135").Value = 496.61905
$ws.Range("K135").Value = 4469.57145
$ws.Range("M135").Value = -1934.57145
$ws.Range("H138").Value = 3051.825
$ws.Range("I138").Value = 1997.0526
$ws.Range("J138").Value = 4006.1428
$ws.Range("K138").Value = 5991.1578
$ws.Range("L138").Value = 12018.4284
$ws.Range("M138").Value = -851.1578
$ws.Range("N138").Value = -22298.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10553.571
$ws.Range("I32").Value = 6214.5894
$ws.Range("J32").Value = 22124.191
$ws.Range("K32").Value = 6214.5894
$ws.Range("L32").Value = 22124.191
$ws.Range("M32").Value = -5927.5894
$ws.Range("N32").Value = -22698.191
$ws.Range("H53").Value = 2039
$ws.Range("I53").Value = 2039
$ws.Range("K53").Value = 2039
$ws.Range("M53").Value = -1357
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H74").Value = 21126.75
$ws.Range("I74").Value = 1966.8536
$ws.Range("K74").Value = 1966.8536
$ws.Range("M74").Value = -1092.8536
$ws.Range("H77").Value = 21126.75
$ws.Range("I77").Value = 1966.8536
$ws.Range("K77").Value = 9834.268
$ws.Range("M77").Value = -5466.268
$ws.Range("H122").Value = 1314243.1
$ws.Range("I122").Value = 1423830
$ws.Range("K122").Value = 4271490
$ws.Range("M122").Value = -4269040
$ws.Range("H139").Value = 89998.5
$ws.Range("J139").Value = 89998.5
$ws.Range("L139").Value = 89998.5
$ws.Range("N139").Value = -100278.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 14196.223
$ws.Range("J81").Value = 14196.223
$ws.Range("L81").Value = 14196.223
$ws.Range("N81").Value = -16318.223
$ws.Range("H84").Value = 14196.223
$ws.Range("J84").Value = 14196.223
$ws.Range("L84").Value = 42588.669
$ws.Range("N84").Value = -53196.669
$ws.Range("H99").Value = 8931397
$ws.Range("I99").Value = 14288006
$ws.Range("J99").Value = 3716.3333
$ws.Range("K99").Value = 14288006
$ws.Range("L99").Value = 3716.3333
$ws.Range("M99").Value = -14286508
$ws.Range("N99").Value = -6712.3333
$ws.Range("H134").Value = 6563.1665
$ws.Range("I134").Value = 1840.3889
$ws.Range("J134").Value = 20731.5
$ws.Range("K134").Value = 5521.1667
$ws.Range("L134").Value = 62194.5
$ws.Range("M134").Value = -2986.1667
$ws.Range("N134").Value = -67264.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17234.209
$ws.Range("I31").Value = 1971.8077
$ws.Range("J31").Value = 26912.805
$ws.Range("K31").Value = 1971.8077
$ws.Range("L31").Value = 26912.805
$ws.Range("M31").Value = -1676.8077
$ws.Range("N31").Value = -27502.805
$ws.Range("H34").Value = 17234.209
$ws.Range("I34").Value = 1971.8077
$ws.Range("J34").Value = 26912.805
$ws.Range("K34").Value = 1971.8077
$ws.Range("L34").Value = 26912.805
$ws.Range("M34").Value = -1769.8077
$ws.Range("N34").Value = -27316.805
$ws.Range("H141").Value = 74896.625
$ws.Range("J141").Value = 83714.14
$ws.Range("L141").Value = 83714.14
$ws.Range("N141").Value = -94074.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 15117806
$ws.Range("I4").Value = 17961274
$ws.Range("J4").Value = 10243290
$ws.Range("K4").Value = 53883822
$ws.Range("L4").Value = 30729870
$ws.Range("M4").Value = -53883710
$ws.Range("N4").Value = -30730094
$ws.Range("H11").Value = 8522.615
$ws.Range("J11").Value = 12571.429
$ws.Range("L11").Value = 37714.287
$ws.Range("N11").Value = -37994.287
$ws.Range("H12").Value = 32502.928
$ws.Range("I12").Value = 68418.69500000001
$ws.Range("J12").Value = 1375.9333
$ws.Range("K12").Value = 205256.085
$ws.Range("L12").Value = 4127.7999
$ws.Range("M12").Value = -205083.085
$ws.Range("N12").Value = -4473.7999
$ws.Range("H26").Value = 314.0625
$ws.Range("I26").Value = 156.07692
$ws.Range("J26").Value = 998.6667
$ws.Range("K26").Value = 468.23076
$ws.Range("L26").Value = 2996.0001
$ws.Range("M26").Value = -180.23076
$ws.Range("N26").Value = -3572.0001
$ws.Range("H114").Value = 501834.6
$ws.Range("I114").Value = 533
$ws.Range("J114").Value = 627160
$ws.Range("K114").Value = 1599
$ws.Range("L114").Value = 1881480
$ws.Range("M114").Value = 1655
$ws.Range("N114").Value = -1887988
$ws.Range("H122").Value = 1264.625
$ws.Range("I122").Value = 1315.6666
$ws.Range("J122").Value = 1111.5
$ws.Range("K122").Value = 11840.9994
$ws.Range("L122").Value = 10003.5
$ws.Range("M122").Value = -9390.999400000001
$ws.Range("N122").Value = -14903.5
$ws.Range("H131").Value = 15436227
$ws.Range("J131").Value = 22227960
$ws.Range("L131").Value = 66683880
$ws.Range("N131").Value = -66693960
$ws.Range("H132").Value = 1268.2069
$ws.Range("I132").Value = 1356.8334
$ws.Range("J132").Value = 1205.6471
$ws.Range("K132").Value = 12211.5006
$ws.Range("L132").Value = 10850.8239
$ws.Range("M132").Value = -9681.500599999999
$ws.Range("N132").Value = -15910.8239

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 5797.2104
$ws.Range("I107").Value = 9535.182000000001
$ws.Range("K107").Value = 9535.182000000001
$ws.Range("M107").Value = -7615.182000000001
$ws.Range("H132").Value = 2409.4795
$ws.Range("I132").Value = 2125.4038
$ws.Range("J132").Value = 3112.9048
$ws.Range("K132").Value = 6376.2114
$ws.Range("L132").Value = 9338.714399999999
$ws.Range("M132").Value = -3846.2114
$ws.Range("N132").Value = -14398.7144
$ws.Range("H136").Value = 26334.143
$ws.Range("J136").Value = 26334.143
$ws.Range("L136").Value = 79002.429
$ws.Range("N136").Value = -84102.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6342.7393
$ws.Range("J46").Value = 6590.136
$ws.Range("L46").Value = 6590.136
$ws.Range("N46").Value = -6966.136
$ws.Range("H55").Value = 2613.8667
$ws.Range("I55").Value = 4242.25
$ws.Range("J55").Value = 2021.7273
$ws.Range("K55").Value = 4242.25
$ws.Range("L55").Value = 2021.7273
$ws.Range("M55").Value = -4069.25
$ws.Range("N55").Value = -2367.7273
$ws.Range("H132").Value = 10110.692
$ws.Range("I132").Value = 10429.521
$ws.Range("K132").Value = 31288.563
$ws.Range("M132").Value = -28758.563
$ws.Range("H136").Value = 39890.555
$ws.Range("I136").Value = 51957.855
$ws.Range("J136").Value = 6906.6
$ws.Range("K136").Value = 155873.565
$ws.Range("L136").Value = 20719.8
$ws.Range("M136").Value = -153323.565
$ws.Range("N136").Value = -25819.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 9964.166999999999
$ws.Range("J74").Value = 9964.166999999999
$ws.Range("L74").Value = 9964.166999999999
$ws.Range("N74").Value = -11836.167
$ws.Range("H77").Value = 9964.166999999999
$ws.Range("J77").Value = 9964.166999999999
$ws.Range("L77").Value = 29892.501
$ws.Range("N77").Value = -39252.501
$ws.Range("H107").Value = 38462516
$ws.Range("I107").Value = 90909784
$ws.Range("K107").Value = 272729352
$ws.Range("M107").Value = -272727432
$ws.Range("H113").Value = 848.0909
$ws.Range("I113").Value = 313.08334
$ws.Range("J113").Value = 1490.1
$ws.Range("K113").Value = 939.2500200000001
$ws.Range("L113").Value = 4470.299999999999
$ws.Range("M113").Value = 1230.74998
$ws.Range("N113").Value = -8810.299999999999
$ws.Range("H122").Value = 3734.36
$ws.Range("I122").Value = 2175.8948
$ws.Range("K122").Value = 6527.6844
$ws.Range("M122").Value = -4077.6844
$ws.Range("H136").Value = 2470.1428
$ws.Range("I136").Value = 2344.725
$ws.Range("K136").Value = 7034.174999999999
$ws.Range("M136").Value = -4484.174999999999
